# final sprint 6 changes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BURNDOWN CHART")

# Update the "Goal or target" starting value (E9). The rest of column E
# (E10:E29) are formulas referencing $E$9 and will recalculate automatically.
$ws.Range("E9").Value = 25

# Update the "Done or actual" values (F9:F28) to the new sprint-6 figures.
$ws.Range("F9").Value = 25
$ws.Range("F10").Value = 25
$ws.Range("F11").Value = 25
$ws.Range("F12").Value = 25
$ws.Range("F13").Value = 25
$ws.Range("F14").Value = 25
$ws.Range("F15").Value = 25
$ws.Range("F16").Value = 25
$ws.Range("F17").Value = 25
$ws.Range("F18").Value = 25
$ws.Range("F19").Value = 25
$ws.Range("F20").Value = 25
$ws.Range("F21").Value = 25
$ws.Range("F22").Value = 25
$ws.Range("F23").Value = 25
$ws.Range("F24").Value = 25
$ws.Range("F25").Value = 25
$ws.Range("F26").Value = 12
$ws.Range("F27").Value = 5
$ws.Range("F28").Value = 5

# Update the chart title text for sprint 6.
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$chart.ChartTitle.Text = "Withdrive Sprint 6 Burndown chart"

# Move the active selection as recorded when the file was last saved.
$ws.Activate()
$ws.Range("M40").Select()
